# "added i2c to ext connector"
#
# New I2C-related parts show up in the BOM:
#   - R56 joins the existing R29, R43 (0 ohm) "Fitted" group on the Top layer
#     (qty 2 -> 3)
#   - A new "Not Fitted" Top-layer row for R57 (1k) is inserted right after the
#     other Top-layer rows (i.e. before the Bottom-layer rows begin)
#   - On the Bottom layer: R40 leaves the R34/R35/.../R52 (1k) group and joins
#     R4 (115k) together with the brand-new R58 (qty 1 -> 3)
#   - R54 and R55 join R7 (8.2k) on the Bottom layer (qty 1 -> 3)
#   - The R36, R38, R46, R47, R51, R53 (0 ohm) Bottom-layer row moves to the
#     very end of the table; it now also covers the new R56 slot, which is
#     itself "Not Fitted" while the rest of the group stays "Fitted"
#     (qty 6 -> 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($rowIndex, $values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
    }
}

function Copy-RowStyle($destRow, $srcRow) {
    $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, 10)).Style = `
        $ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, 10)).Style
    $ws.Cells.Item($destRow, 5).Style = $ws.Cells.Item($srcRow, 5).Style
}

# --- Row 24 (Top layer): R29, R43 -> R29, R43, R56, quantity 2 -> 3 ---
$ws.Cells.Item(24, 3).Value = "R29, R43, R56"
$ws.Cells.Item(24, 5).Value = 3

# --- Insert a brand-new row 34 (Top layer, R57 "Not Fitted", qty 0) ---
# Inserting at row 34 shifts the previous rows 34-60 (the whole "Bottom" block)
# down to rows 35-61, keeping their content and formatting intact.
$ws.Rows(34).Insert()
Set-RowValues 34 @("Top", "Not Fitted", "R57", "1k", 0, "CRCW06031K00FKEA", "RESC1608N", "100mW,75V", "1469740", "Farnell")
Copy-RowStyle 34 33

# --- Bottom layer (rows shifted by the insert above) ---

# Row 51: R34, R35, R40, R41, R44, R45, R52 (1k, qty 7) -> R40 leaves the group
$ws.Cells.Item(51, 3).Value = "R34, R35, R41, R44, R45, R52"
$ws.Cells.Item(51, 5).Value = 6

# Row 52 currently holds the R36, R38, R46, R47, R51, R53 (0 ohm) group; that
# whole row relocates to the end of the table, so remove it here (rows 53-61
# shift up to 52-60).
$ws.Rows(52).Delete()

# Row 52 is now R4 (115k, qty 1) -> R4, R40, R58 (qty 3)
$ws.Cells.Item(52, 3).Value = "R4, R40, R58"
$ws.Cells.Item(52, 5).Value = 3

# Row 54 is now R7 (8.2k, qty 1) -> R7, R54, R55 (qty 3)
$ws.Cells.Item(54, 3).Value = "R7, R54, R55"
$ws.Cells.Item(54, 5).Value = 3

# Append the relocated R36, R38, R46, R47, R51, R53 (0 ohm) row as the new last
# row (61), updated to reflect the new "Not Fitted" R56 slot it now covers.
$ws.Rows(61).Insert()
Set-RowValues 61 @("Bottom", "Not Fitted, Fitted, Fitted, Fitted, Fitted, Fitted", "R36, R38, R46, R47, R51, R53", "0", 5, "CRCW06030000Z0EA", "RESC1608N", "100mW,75V", "1469739", "Farnell")
Copy-RowStyle 61 60

$ws.Range("A1:J61").Columns.AutoFit() | Out-Null
